$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.339.91'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '3.538.17'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.29'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.00'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '3.536.12'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -0.83%  '
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.380'
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").Value = '4.133.51'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.69'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("D17").Value = '3.531.80'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '64.319.00'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.82'
$ws.Range("E19").Value = '  -3.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.10'
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.65'
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.61'
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.577'
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").Value = '3.679.20'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.20'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("E29").Value = '  -3.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.40'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").Value = '3.543.35'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.68'
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.40'
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.95'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.35'
$ws.Range("E40").Value = '  -4.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0791'
$ws.Range("E41").Value = '  -2.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.817'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.47'
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.13'
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.43'
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("D49").Value = '2.474.05'
$ws.Range("E49").Value = '  +2.06%  '
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.908'
$ws.Range("E51").Value = '  -0.43%  '
